$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.201.96'
$ws.Range("E2").Value = '  +1.40%  '
$ws.Range("D3").Value = '1.609.25'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.46'
$ws.Range("E5").Value = '  +1.45%  '
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.482'
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("E8").Value = '  +1.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0619'
$ws.Range("E9").Value = '  +1.49%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.24'
$ws.Range("E10").Value = '  +1.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0797'
$ws.Range("E11").Value = '  +1.68%  '
$ws.Range("D12").Value = '1.834.14'
$ws.Range("E12").Value = '  +0.67%  '
$ws.Range("D13").Value = '1.611.78'
$ws.Range("E13").Value = '  +2.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.02'
$ws.Range("E14").Value = '  -0.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.511'
$ws.Range("E15").Value = '  +0.35%  '
$ws.Range("D16").Value = '26.209.95'
$ws.Range("E16").Value = '  +1.48%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.77'
$ws.Range("E17").Value = '  +0.45%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '0.0₃0729'
$ws.Range("E18").Value = '  +1.67%  '
$ws.Range("E19").Value = '  -0.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '199.98'
$ws.Range("E20").Value = '  +5.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.25'
$ws.Range("E21").Value = '  +1.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.43'
$ws.Range("E22").Value = '  +1.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.02'
$ws.Range("E23").Value = '  +1.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.132'
$ws.Range("E24").Value = '  +3.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.24'
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("E26").Value = '  +2.15%  '
$ws.Range("E27").Value = '  -0.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.20'
$ws.Range("E28").Value = '  +1.72%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.50'
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("E30").Value = '  -0.98%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0475'
$ws.Range("E31").Value = '  +1.81%  '
$ws.Range("E32").Value = '  +1.79%  '
$ws.Range("E33").Value = '  +1.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.50'
$ws.Range("E34").Value = '  +2.21%  '
$ws.Range("E35").Value = '  -1.54%  '
$ws.Range("D36").Value = '1.107.56'
$ws.Range("E36").Value = '  +0.81%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.36'
$ws.Range("E37").Value = '  -0.20%  '
$ws.Range("E38").Value = '  +0.96%  '
$ws.Range("E39").Value = '  -0.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.503'
$ws.Range("E40").Value = '  +1.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.785'
$ws.Range("E41").Value = '  -0.87%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.790'
$ws.Range("E42").Value = '  +6.41%  '
$ws.Range("D43").Value = '1.746.15'
$ws.Range("E43").Value = '  +0.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.15'
$ws.Range("E44").Value = '  +1.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '93.07'
$ws.Range("E45").Value = '  -2.66%  '
$ws.Range("D46").Value = '0.0₆0106'
$ws.Range("E46").Value = '  -5.73%  '
$ws.Range("E47").Value = '  +7.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.69'
$ws.Range("E48").Value = '  +0.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0510'
$ws.Range("E49").Value = '  -0.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.410'
$ws.Range("E50").Value = '  +0.35%  '
$ws.Range("E51").Value = '  -0.15%  '